$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.484.76"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.429.62"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "233.71"
$ws.Range("E5").Value = "  -1.23%  "

$ws.Range("D6").Value = "621.31"
$ws.Range("E6").Value = "  -2.71%  "

$ws.Range("E7").Value = "  -2.22%  "

$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "0.976"
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").Value = "3.428.34"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "42.99"
$ws.Range("E12").Value = "  +4.11%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "6.28"
$ws.Range("E14").Value = "  +1.99%  "

$ws.Range("D15").Value = "93.295.00"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "4.067.85"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").Value = "8.21"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").Value = "3.428.92"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "18.12"
$ws.Range("E20").Value = "  +4.51%  "

$ws.Range("D21").Value = "11.68"
$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").Value = "502.76"
$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("D23").Value = "3.38"
$ws.Range("E23").Value = "  +5.06%  "

$ws.Range("D24").Value = "0.451"
$ws.Range("E24").Value = "  -3.01%  "

$ws.Range("D25").Value = "6.64"
$ws.Range("E25").Value = "  +3.05%  "

$ws.Range("E26").Value = "  -2.83%  "

$ws.Range("D27").Value = "95.01"
$ws.Range("E27").Value = "  +4.56%  "

$ws.Range("D28").Value = "11.98"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("D29").Value = "3.609.54"
$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("D30").Value = "11.45"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  +2.21%  "

$ws.Range("D33").Value = "2.75"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").Value = "0.994"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").Value = "30.03"
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("D37").Value = "0.552"
$ws.Range("E37").Value = "  +1.47%  "

$ws.Range("D38").Value = "557.00"
$ws.Range("E38").Value = "  +3.79%  "

$ws.Range("D39").Value = "7.50"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.149"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.915"
$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").Value = "3.69"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "5.51"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").Value = "0.0411"
$ws.Range("E48").Value = "  +2.32%  "

$ws.Range("E49").Value = "  -1.16%  "

$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("D51").Value = "8.10"
$ws.Range("E51").Value = "  +1.23%  "
